$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mail: "Close" label row, inserted before the existing row 92 -----------
# Shifts every row from the old 92 downward by one (old 92..106 -> 93..107),
# keeping the blank separator rows that already existed in the sheet.
$ws.Rows.Item(92).Insert()

$ws.Range("A92").Value = "lang_close"
$ws.Range("B92").Value = "Đóng"
$ws.Range("C92").Value = "Close"

# --- Mail: "Confirm delete mail ?" row, appended after the old last row -----
# The old last row (106) is now 107; leave row 108 blank (matching the
# existing single blank-row separator convention) and write the new data on
# row 109.
$ws.Range("A109").Value = "lang_confirm_delete_mail"
$ws.Range("B109").Value = "Xác nhận xóa thư ?"
$ws.Range("C109").Value = "Confirm delete Mail ?"

# --- Match the author's final view state (best effort) ---------------------
[void]$ws.Range("C109").Select()
$excel.ActiveWindow.ScrollRow = 102
